# Adds the new donation report rows (16-26) captured after the test fixes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Samuel", "sam.abreugouveia@gmail.com", "(81) 99899-1917", "Plástico", "Gostaria de doar plástico", "26/11/2024 17:30"),
    @("Samuel", "sam.abreugouveia@gmail.com", "(81) 99899-1917", "Plástico", "Eu gostaria de doar plástico", "26/11/2024 17:30"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 17:44"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 17:45"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 17:46"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 17:49"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 17:50"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 17:50"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 18:04"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 18:06"),
    @("Doador", "doador@gmail.com", "(81) 99999-9999", "Resíduos têxteis", "Gostaria de fazer uma doação de resíduos têxteis para a confecção de novos brinquedos.", "26/11/2024 18:22")
)

$startRow = 16
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
